$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "day" column (AE) after the current last column (AD) with the
# attendance code for each player, mirroring the formatting of column AD.
# ---------------------------------------------------------------------------

# Correct a previous data-entry mistake: row 19 (Jeremie Laurent) should be
# marked "B" (Blessure) on the AD column day, not "RH".
$ws.Cells.Item(19, 30).Value = "B"

# New date header for the added day (19 Aug 2025 -> 20 Aug 2025).
$ws.Cells.Item(1, 31).Value = 45889

# Attendance codes for the new day, one per player row (rows 2-27).
$attendance = @{
    2  = "P"
    3  = "M"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    12 = "P"
    13 = "A"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "B"
    18 = "P"
    19 = "B"
    20 = "P"
    21 = "P"
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "RH"
}

foreach ($row in $attendance.Keys) {
    $ws.Cells.Item($row, 31).Value = $attendance[$row]
}

# Copy number formats / styles from column AD (rows 1-27) onto the new
# column AE so the new date header and data cells look like the rest of the
# table. This is done *after* the values are written so the formula cells
# (B:I) have already recalculated against the new data.
$ws.Range("AD1:AD27").Copy()
$ws.Range("AE1:AE27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active cell / selection to match where the user left off.
$ws.Activate()
$ws.Range("AG23").Select()
